$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"1"
$ws.Range("F2").Value = [double]"0.3333333333333333"
$ws.Range("G2").Value = [double]"0.4460406666666667"
$ws.Range("H2").Value = [double]"1.338122"
$ws.Range("I2").Value = [double]"0.001628842357811546"
$ws.Range("J2").Value = [double]"0.001628842357811545"
$ws.Range("M2").Value = [double]"0.243056"
$ws.Range("N2").Value = [double]"0.729168"
$ws.Range("O2").Value = [double]"0.002199620488481675"
$ws.Range("P2").Value = [double]"0.002199620488481675"
$ws.Range("Q2").Value = [double]"0.1084128602773334"
$ws.Range("R2").Value = [double]"0.9757157424960001"
$ws.Range("S2").Value = [double]"3.582835022749074E-06"
$ws.Range("T2").Value = [double]"3.582835022749074E-06"
$ws.Range("E3").Value = [double]"1"
$ws.Range("F3").Value = [double]"0.3333333333333333"
$ws.Range("G3").Value = [double]"0.4460406666666667"
$ws.Range("H3").Value = [double]"1.338122"
$ws.Range("I3").Value = [double]"0.001628842357811546"
$ws.Range("J3").Value = [double]"0.001628842357811545"
$ws.Range("M3").Value = [double]"70.95253000000001"
$ws.Range("N3").Value = [double]"212.85759"
$ws.Range("O3").Value = [double]"0.6421097964979703"
$ws.Range("P3").Value = [double]"0.6421097964979703"
$ws.Range("Q3").Value = [double]"31.64771378288667"
$ws.Range("R3").Value = [double]"284.82942404598"
$ws.Range("S3").Value = [double]"0.001045895634901646"
$ws.Range("T3").Value = [double]"0.001045895634901645"
$ws.Range("E4").Value = [double]"1"
$ws.Range("F4").Value = [double]"0.3333333333333333"
$ws.Range("G4").Value = [double]"0.4460406666666667"
$ws.Range("H4").Value = [double]"1.338122"
$ws.Range("I4").Value = [double]"0.001628842357811546"
$ws.Range("J4").Value = [double]"0.001628842357811545"
$ws.Range("M4").Value = [double]"0.04794200000000001"
$ws.Range("N4").Value = [double]"0.143826"
$ws.Range("O4").Value = [double]"0.0004338679376719292"
$ws.Range("P4").Value = [double]"0.0004338679376719292"
$ws.Range("Q4").Value = [double]"0.02138408164133334"
$ws.Range("R4").Value = [double]"0.192456734772"
$ws.Range("S4").Value = [double]"7.067024745763779E-07"
$ws.Range("T4").Value = [double]"7.067024745763778E-07"
$ws.Range("E5").Value = [double]"1"
$ws.Range("F5").Value = [double]"0.3333333333333333"
$ws.Range("G5").Value = [double]"0.4460406666666667"
$ws.Range("H5").Value = [double]"1.338122"
$ws.Range("I5").Value = [double]"0.001628842357811546"
$ws.Range("J5").Value = [double]"0.001628842357811545"
$ws.Range("M5").Value = [double]"39.25553366666666"
$ws.Range("N5").Value = [double]"117.766601"
$ws.Range("O5").Value = [double]"0.3552567150758761"
$ws.Range("P5").Value = [double]"0.3552567150758761"
$ws.Range("Q5").Value = [double]"17.50956440703578"
$ws.Range("R5").Value = [double]"157.586079663322"
$ws.Range("S5").Value = [double]"0.0005786571854125744"
$ws.Range("T5").Value = [double]"0.0005786571854125744"
$ws.Range("G6").Value = [double]"145.8660203333333"
$ws.Range("H6").Value = [double]"437.598061"
$ws.Range("I6").Value = [double]"0.5326706066061244"
$ws.Range("J6").Value = [double]"0.5326706066061244"
$ws.Range("M6").Value = [double]"0.243056"
$ws.Range("N6").Value = [double]"0.729168"
$ws.Range("O6").Value = [double]"0.002199620488481675"
$ws.Range("P6").Value = [double]"0.002199620488481675"
$ws.Range("Q6").Value = [double]"35.45361143813867"
$ws.Range("R6").Value = [double]"319.082502943248"
$ws.Range("S6").Value = [double]"0.001171673179902793"
$ws.Range("T6").Value = [double]"0.001171673179902793"
$ws.Range("G7").Value = [double]"145.8660203333333"
$ws.Range("H7").Value = [double]"437.598061"
$ws.Range("I7").Value = [double]"0.5326706066061244"
$ws.Range("J7").Value = [double]"0.5326706066061244"
$ws.Range("M7").Value = [double]"70.95253000000001"
$ws.Range("N7").Value = [double]"212.85759"
$ws.Range("O7").Value = [double]"0.6421097964979703"
$ws.Range("P7").Value = [double]"0.6421097964979703"
$ws.Range("Q7").Value = [double]"10349.56318368144"
$ws.Range("R7").Value = [double]"93146.06865313298"
$ws.Range("S7").Value = [double]"0.3420330148083089"
$ws.Range("T7").Value = [double]"0.3420330148083089"
$ws.Range("G8").Value = [double]"145.8660203333333"
$ws.Range("H8").Value = [double]"437.598061"
$ws.Range("I8").Value = [double]"0.5326706066061244"
$ws.Range("J8").Value = [double]"0.5326706066061244"
$ws.Range("M8").Value = [double]"0.04794200000000001"
$ws.Range("N8").Value = [double]"0.143826"
$ws.Range("O8").Value = [double]"0.0004338679376719292"
$ws.Range("P8").Value = [double]"0.0004338679376719292"
$ws.Range("Q8").Value = [double]"6.993108746820667"
$ws.Range("R8").Value = [double]"62.937978721386"
$ws.Range("S8").Value = [double]"0.0002311086975466547"
$ws.Range("T8").Value = [double]"0.0002311086975466547"
$ws.Range("G9").Value = [double]"145.8660203333333"
$ws.Range("H9").Value = [double]"437.598061"
$ws.Range("I9").Value = [double]"0.5326706066061244"
$ws.Range("J9").Value = [double]"0.5326706066061244"
$ws.Range("M9").Value = [double]"39.25553366666666"
$ws.Range("N9").Value = [double]"117.766601"
$ws.Range("O9").Value = [double]"0.3552567150758761"
$ws.Range("P9").Value = [double]"0.3552567150758761"
$ws.Range("Q9").Value = [double]"5726.04847201785"
$ws.Range("R9").Value = [double]"51534.43624816066"
$ws.Range("S9").Value = [double]"0.189234809920366"
$ws.Range("T9").Value = [double]"0.189234809920366"
$ws.Range("E10").Value = [double]"2"
$ws.Range("F10").Value = [double]"0.6666666666666666"
$ws.Range("G10").Value = [double]"0.150912"
$ws.Range("H10").Value = [double]"0.452736"
$ws.Range("I10").Value = [double]"0.0005510974139175409"
$ws.Range("J10").Value = [double]"0.0005510974139175409"
$ws.Range("M10").Value = [double]"0.243056"
$ws.Range("N10").Value = [double]"0.729168"
$ws.Range("O10").Value = [double]"0.002199620488481675"
$ws.Range("P10").Value = [double]"0.002199620488481675"
$ws.Range("Q10").Value = [double]"0.036680067072"
$ws.Range("R10").Value = [double]"0.330120603648"
$ws.Range("S10").Value = [double]"1.212205162802289E-06"
$ws.Range("T10").Value = [double]"1.212205162802289E-06"
$ws.Range("E11").Value = [double]"2"
$ws.Range("F11").Value = [double]"0.6666666666666666"
$ws.Range("G11").Value = [double]"0.150912"
$ws.Range("H11").Value = [double]"0.452736"
$ws.Range("I11").Value = [double]"0.0005510974139175409"
$ws.Range("J11").Value = [double]"0.0005510974139175409"
$ws.Range("M11").Value = [double]"70.95253000000001"
$ws.Range("N11").Value = [double]"212.85759"
$ws.Range("O11").Value = [double]"0.6421097964979703"
$ws.Range("P11").Value = [double]"0.6421097964979703"
$ws.Range("Q11").Value = [double]"10.70758820736"
$ws.Range("R11").Value = [double]"96.36829386623999"
$ws.Range("S11").Value = [double]"0.0003538650483011499"
$ws.Range("T11").Value = [double]"0.0003538650483011499"
$ws.Range("E12").Value = [double]"2"
$ws.Range("F12").Value = [double]"0.6666666666666666"
$ws.Range("G12").Value = [double]"0.150912"
$ws.Range("H12").Value = [double]"0.452736"
$ws.Range("I12").Value = [double]"0.0005510974139175409"
$ws.Range("J12").Value = [double]"0.0005510974139175409"
$ws.Range("M12").Value = [double]"0.04794200000000001"
$ws.Range("N12").Value = [double]"0.143826"
$ws.Range("O12").Value = [double]"0.0004338679376719292"
$ws.Range("P12").Value = [double]"0.0004338679376719292"
$ws.Range("Q12").Value = [double]"0.007235023104000001"
$ws.Range("R12").Value = [double]"0.06511520793599999"
$ws.Range("S12").Value = [double]"2.391034984327371E-07"
$ws.Range("T12").Value = [double]"2.391034984327371E-07"
$ws.Range("E13").Value = [double]"2"
$ws.Range("F13").Value = [double]"0.6666666666666666"
$ws.Range("G13").Value = [double]"0.150912"
$ws.Range("H13").Value = [double]"0.452736"
$ws.Range("I13").Value = [double]"0.0005510974139175409"
$ws.Range("J13").Value = [double]"0.0005510974139175409"
$ws.Range("M13").Value = [double]"39.25553366666666"
$ws.Range("N13").Value = [double]"117.766601"
$ws.Range("O13").Value = [double]"0.3552567150758761"
$ws.Range("P13").Value = [double]"0.3552567150758761"
$ws.Range("Q13").Value = [double]"5.924131096704"
$ws.Range("R13").Value = [double]"53.31717987033599"
$ws.Range("S13").Value = [double]"0.000195781056955156"
$ws.Range("T13").Value = [double]"0.000195781056955156"
$ws.Range("G14").Value = [double]"127.376091"
$ws.Range("H14").Value = [double]"382.128273"
$ws.Range("I14").Value = [double]"0.4651494536221465"
$ws.Range("J14").Value = [double]"0.4651494536221465"
$ws.Range("M14").Value = [double]"0.243056"
$ws.Range("N14").Value = [double]"0.729168"
$ws.Range("O14").Value = [double]"0.002199620488481675"
$ws.Range("P14").Value = [double]"0.002199620488481675"
$ws.Range("Q14").Value = [double]"30.95952317409601"
$ws.Range("R14").Value = [double]"278.635708566864"
$ws.Range("S14").Value = [double]"0.00102315226839333"
$ws.Range("T14").Value = [double]"0.00102315226839333"
$ws.Range("G15").Value = [double]"127.376091"
$ws.Range("H15").Value = [double]"382.128273"
$ws.Range("I15").Value = [double]"0.4651494536221465"
$ws.Range("J15").Value = [double]"0.4651494536221465"
$ws.Range("M15").Value = [double]"70.95253000000001"
$ws.Range("N15").Value = [double]"212.85759"
$ws.Range("O15").Value = [double]"0.6421097964979703"
$ws.Range("P15").Value = [double]"0.6421097964979703"
$ws.Range("Q15").Value = [double]"9037.655917960232"
$ws.Range("R15").Value = [double]"81338.90326164209"
$ws.Range("S15").Value = [double]"0.2986770210064586"
$ws.Range("T15").Value = [double]"0.2986770210064585"
$ws.Range("G16").Value = [double]"127.376091"
$ws.Range("H16").Value = [double]"382.128273"
$ws.Range("I16").Value = [double]"0.4651494536221465"
$ws.Range("J16").Value = [double]"0.4651494536221465"
$ws.Range("M16").Value = [double]"0.04794200000000001"
$ws.Range("N16").Value = [double]"0.143826"
$ws.Range("O16").Value = [double]"0.0004338679376719292"
$ws.Range("P16").Value = [double]"0.0004338679376719292"
$ws.Range("Q16").Value = [double]"6.106664554722001"
$ws.Range("R16").Value = [double]"54.95998099249801"
$ws.Range("S16").Value = [double]"0.0002018134341522654"
$ws.Range("T16").Value = [double]"0.0002018134341522654"
$ws.Range("G17").Value = [double]"127.376091"
$ws.Range("H17").Value = [double]"382.128273"
$ws.Range("I17").Value = [double]"0.4651494536221465"
$ws.Range("J17").Value = [double]"0.4651494536221465"
$ws.Range("M17").Value = [double]"39.25553366666666"
$ws.Range("N17").Value = [double]"117.766601"
$ws.Range("O17").Value = [double]"0.3552567150758761"
$ws.Range("P17").Value = [double]"0.3552567150758761"
$ws.Range("Q17").Value = [double]"5000.216428578898"
$ws.Range("R17").Value = [double]"45001.94785721008"
$ws.Range("S17").Value = [double]"0.1652474669131423"
$ws.Range("T17").Value = [double]"0.1652474669131423"
